# "parallel execution and test fixes"
# Update the "Data" sheet: drop the old "Brand" label in A2, and add a
# new Watch / "Jewelry & Watches" row of data in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# A2 ("Brand") is no longer used - remove it entirely.
$ws.Range("A2").ClearContents()

# New third row of data: Watch / Jewelry & Watches.
$ws.Range("C3").Value = "Watch"
$ws.Range("C4").Value = "Jewelry & Watches"

# Column B needs to be wider now that it holds longer category text.
$ws.Columns.Item(2).ColumnWidth = 25.3

# Selection moves to the full data block.
[void]$ws.Range("B2:C4").Select()
